# "remove none as a climate choice"
#
# The "choices" sheet has a `climate` choice list starting at row 42 whose
# first entry is name="none" / label="None". That row is removed entirely
# (every row below shifts up by one), and the `survey` sheet's
# select_multiple-climate question (row 12) had "none" set as its XLSForm
# `default` (column H) - that default is cleared since the choice no
# longer exists.

$wb = $excel.ActiveWorkbook

# --- choices sheet: drop the "climate / none / None" row -------------------
$choices = $wb.Worksheets.Item("choices")
$choices.Activate()
$choices.Rows.Item(42).Delete()
$choices.Rows.Item(42).Select()

# --- survey sheet: clear the now-invalid "none" default --------------------
$survey = $wb.Worksheets.Item("survey")
$survey.Activate()
$survey.Range("H12").ClearContents()
$survey.Range("H12").Select()
